# Updates the cryptos list values/percentages on Sheet1 to reflect the
# refreshed data from the GitHub Actions run. For "Price" column (D) values
# that look like plain numbers, the cell is forced to Text format first so
# Excel does not silently convert the string into a numeric value (which
# would strip meaningful trailing zeros / punctuation such as "1.00" -> 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.922.58'
$ws.Range("E2").Value = '  +1.10%  '
$ws.Range("D3").Value = '2.288.92'
$ws.Range("E3").Value = '  +2.37%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '252.43'
$ws.Range("E5").Value = '  +0.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.642'
$ws.Range("E6").Value = '  +1.54%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '73.88'
$ws.Range("E7").Value = '  +6.52%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.648'
$ws.Range("E9").Value = '  +2.42%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.94'
$ws.Range("E10").Value = '  -4.39%  '
$ws.Range("E11").Value = '  +3.07%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '59.06'
$ws.Range("E12").Value = '  -0.43%  '
$ws.Range("E13").Value = '  +2.29%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.106'
$ws.Range("E14").Value = '  +2.71%  '
$ws.Range("D15").Value = '2.631.58'
$ws.Range("E15").Value = '  +2.40%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.34'
$ws.Range("E16").Value = '  +5.00%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.875'
$ws.Range("E17").Value = '  -1.39%  '
$ws.Range("D18").Value = '2.292.69'
$ws.Range("E18").Value = '  +1.72%  '
$ws.Range("D19").Value = '42.829.02'
$ws.Range("E19").Value = '  +1.29%  '
$ws.Range("E20").Value = '  +3.69%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.31'
$ws.Range("E21").Value = '  +1.62%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '72.82'
$ws.Range("E22").Value = '  +0.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.12'
$ws.Range("E23").Value = '  +1.26%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.22'
$ws.Range("E24").Value = '  +6.75%  '
$ws.Range("E25").Value = '  -1.40%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.58'
$ws.Range("E26").Value = '  +0.55%  '
$ws.Range("E27").Value = '  -0.25%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.43'
$ws.Range("E28").Value = '  +0.07%  '
$ws.Range("E29").Value = '  -0.64%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.15'
$ws.Range("E30").Value = '  -2.67%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '166.99'
$ws.Range("E31").Value = '  -0.23%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.07'
$ws.Range("E32").Value = '  +1.21%  '
$ws.Range("E33").Value = '  +6.08%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.128'
$ws.Range("E34").Value = '  +3.78%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0829'
$ws.Range("E35").Value = '  +5.80%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '31.13'
$ws.Range("E36").Value = '  +11.96%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.127'
$ws.Range("E37").Value = '  +2.38%  '
$ws.Range("E38").Value = '  +12.18%  '
$ws.Range("E39").Value = '  +2.24%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0310'
$ws.Range("E40").Value = '  -2.92%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '14.40'
$ws.Range("E41").Value = '  +14.61%  '
$ws.Range("E42").Value = '  +3.68%  '
$ws.Range("E43").Value = '  +3.85%  '
$ws.Range("E44").Value = '  +8.37%  '
$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '9.13'
$ws.Range("E45").Value = '  +4.76%  '
$ws.Range("B46").Value = 'FTXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.88'
$ws.Range("E46").Value = '  -3.67%  '
$ws.Range("B47").Value = 'MultiversX'
$ws.Range("C47").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '60.46'
$ws.Range("E47").Value = '  -4.80%  '
$ws.Range("E48").Value = '  +2.26%  '
$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '101.42'
$ws.Range("E49").Value = '  +7.49%  '
$ws.Range("B50").Value = 'BinanceUSD'
$ws.Range("C50").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.00'
$ws.Range("E50").Value = '  +0.15%  '
$ws.Range("E51").Value = '  +1.21%  '
